# "Cargue un script y subi a bigboo"
#
# 1. The first paragraph ("Hhola" + spell-check markers + " como estas",
#    split across two runs) becomes a single clean run reading
#    "Hhola como estas", with the w:proofErr spell-check markers removed.
# 2. A new blank paragraph is added.
# 3. A new paragraph with the text "Agus agre un scrpto" is added.

$d = $word.ActiveDocument

# --- Paragraph 1: clear the old runs/proofErr markers and retype the text ---
$p1 = $d.Paragraphs(1).Range
$p1.Delete()

$p1 = $d.Paragraphs(1).Range
$p1.InsertAfter("Hhola como estas")

# --- Append the blank paragraph and the new "Agus agre un scrpto" paragraph ---
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParas = "<w:p $wNs/>" + `
            "<w:p $wNs><w:r><w:t>Agus agre un scrpto</w:t></w:r></w:p>"

$tail = $d.Content
$tail.Collapse(0)
$tail.InsertXML($newParas)
